$wb = $excel.ActiveWorkbook

# --- Sheet 1 (treatment): update existing meta-analysis values with refined precision ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B3").Value = 0.5959455751686174
$ws1.Range("C3").Value = 1.013994367340575
$ws1.Range("D3").Value = 1.615694245733231
$ws1.Range("B4").Value = 0.4618510332094891
$ws1.Range("C4").Value = 1.287831479938114
$ws1.Range("D4").Value = 2.529853038381398
$ws1.Range("B5").Value = 0.6941222355568236
$ws1.Range("C5").Value = 1.159082944814026
$ws1.Range("D5").Value = 1.6245483587712
$ws1.Range("B7").Value = 26.32521348497664
$ws1.Range("C7").Value = 35.813664395729
$ws1.Range("D7").Value = 47.62549860435363
$ws1.Range("B8").Value = 0.1999870001396751
$ws1.Range("C8").Value = 0.461383181708975
$ws1.Range("D8").Value = 0.8301648248825331
$ws1.Range("B9").Value = 0.4540789713723684
$ws1.Range("C9").Value = 0.6897018435148856
$ws1.Range("D9").Value = 0.9251511029733821
$ws1.Range("B11").Value = 6.123347660593229
$ws1.Range("C11").Value = 8.775997425171873
$ws1.Range("D11").Value = 12.17825231403138
$ws1.Range("B12").Value = 0.2353306090000331
$ws1.Range("C12").Value = 0.5852165482337847
$ws1.Range("D12").Value = 1.091708863976712
$ws1.Range("B13").Value = 0.4936927331846208
$ws1.Range("C13").Value = 0.7785308463259701
$ws1.Range("D13").Value = 1.063337043478054

# --- Sheet 1 (treatment): add new "Speed meta analysis" block (rows 14-17) ---
$ws1.Range("A14").Value = "Speed meta analysis"
$ws1.Range("A15").Value = "mean (km/day)"
$ws1.Range("B15").Value = 4.112431376685026
$ws1.Range("C15").Value = 4.760985039628134
$ws1.Range("D15").Value = 5.480449112430368
$ws1.Range("A16").Value = "CoV² (RVAR)"
$ws1.Range("B16").Value = 0.03542117434183563
$ws1.Range("C16").Value = 0.0851651924223696
$ws1.Range("D16").Value = 0.1563522747651325
$ws1.Range("A17").Value = "CoV  (RSTD)"
$ws1.Range("B17").Value = 0.191337966285609
$ws1.Range("C17").Value = 0.2966885197628259
$ws1.Range("D17").Value = 0.4019960427239915

# --- Sheet 2 (control): update existing meta-analysis values with refined precision ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = 41.0841764212609
$ws2.Range("C3").Value = 52.24317489057803
$ws2.Range("D3").Value = 65.47938736384387
$ws2.Range("B4").Value = 0.07954459666577422
$ws2.Range("C4").Value = 0.2073836971666101
$ws2.Range("D4").Value = 0.3954039972050888
$ws2.Range("B5").Value = 0.2874496393090814
$ws2.Range("C5").Value = 0.4641343725828301
$ws2.Range("D5").Value = 0.640880346703167
$ws2.Range("B7").Value = 19.1068835724091
$ws2.Range("C7").Value = 22.32140647266363
$ws2.Range("D7").Value = 25.91730203741737
$ws2.Range("B8").Value = 0.03582192564724802
$ws2.Range("C8").Value = 0.08946553009626575
$ws2.Range("D8").Value = 0.1672327581618985
$ws2.Range("B9").Value = 0.192641405648309
$ws2.Range("C9").Value = 0.3044409329938641
$ws2.Range("D9").Value = 0.4162323689434379
$ws2.Range("B11").Value = 5.318410720789306
$ws2.Range("C11").Value = 6.523880458632495
$ws2.Range("D11").Value = 7.916469803675057
$ws2.Range("B12").Value = 0.054173483368058
$ws2.Range("C12").Value = 0.1463465908338701
$ws2.Range("D12").Value = 0.2835157999455275
$ws2.Range("B13").Value = 0.2374859503411489
$ws2.Range("C13").Value = 0.3903335181468861
$ws2.Range("D13").Value = 0.5432918301291978

# --- Sheet 2 (control): add new "Speed meta analysis" block (rows 14-17) ---
$ws2.Range("A14").Value = "Speed meta analysis"
$ws2.Range("A15").Value = "mean (km/day)"
$ws2.Range("B15").Value = 3.633749656974639
$ws2.Range("C15").Value = 4.280444381018343
$ws2.Range("D15").Value = 5.005375851837153
$ws2.Range("A16").Value = "CoV² (RVAR)"
$ws2.Range("B16").Value = 0.02156303552970272
$ws2.Range("C16").Value = 0.07173901326442587
$ws2.Range("D16").Value = 0.151526132229185
$ws2.Range("A17").Value = "CoV  (RSTD)"
$ws2.Range("B17").Value = 0.1508584116135067
$ws2.Range("C17").Value = 0.2751644225042091
$ws2.Range("D17").Value = 0.3999063961861548

Write-Output "done"
